$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lo = $ws.ListObjects.Item(1)
$newCol = $lo.ListColumns.Add()

$ws.Range("D1").Value = "Gracia"
$ws.Range("D2").Value = 0
$ws.Range("D3").Value = 30
$ws.Range("D4").Value = 40

$ws.Range("D3").Font.Name = "Calibri"

$ws.Range("D3").Select() | Out-Null
